$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2025-07-14 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-07-15 Tuesday", 2)

# Update the division problems in the table. Each value is addressed by its
# exact row/column location to avoid ambiguity, since several new values
# collide with other old/new values elsewhere in the table.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "37÷4="
$t.Cell(1,2).Range.Text = "61÷9="
$t.Cell(1,3).Range.Text = "34÷8="
$t.Cell(1,4).Range.Text = "54÷8="
$t.Cell(1,5).Range.Text = "85÷8="

$t.Cell(5,1).Range.Text = "58÷3="
$t.Cell(5,2).Range.Text = "95÷2="
$t.Cell(5,3).Range.Text = "13÷2="
$t.Cell(5,4).Range.Text = "22÷5="
$t.Cell(5,5).Range.Text = "92÷8="

$t.Cell(9,1).Range.Text = "69÷9="
$t.Cell(9,2).Range.Text = "20÷4="
$t.Cell(9,3).Range.Text = "68÷4="
$t.Cell(9,4).Range.Text = "90÷2="
$t.Cell(9,5).Range.Text = "82÷5="

$t.Cell(13,1).Range.Text = "69÷6="
$t.Cell(13,2).Range.Text = "61÷9="
$t.Cell(13,3).Range.Text = "24÷4="
$t.Cell(13,4).Range.Text = "62÷8="
$t.Cell(13,5).Range.Text = "16÷6="

$t.Cell(17,1).Range.Text = "74÷8="
$t.Cell(17,2).Range.Text = "85÷5="
$t.Cell(17,3).Range.Text = "91÷7="
$t.Cell(17,4).Range.Text = "46÷3="
$t.Cell(17,5).Range.Text = "60÷3="
